$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AdduserData")

# Update data row 2 values
$ws.Range("A2").Value = "kala"
$ws.Range("B2").Value = "khatta "
$ws.Range("C2").Value = "jamun"
$ws.Range("D2").Value = 104
$ws.Range("E2").Value = "jamun"

# Move the active selection from D6 to D2
[void]$ws.Activate()
[void]$ws.Range("D2").Select()

# Adjust the workbook window width (best effort; mirrors the author's resize)
$excel.ActiveWindow.Width = 16872
